{"js": "// Bill revision: 00017 -> 00018, due date 2025-06-19 -> 2025-06-18,\n// meter reading Previous 324.0 -> 435.0, Present 435.0 -> 546.0.\n// The document contains two identical copies of the bill, so every\n// target string is searched/replaced in document order (2 hits each).\n\nasync function replaceAll(context, searchText, newText) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Order matters: do the \"435.0 -> 546.0\" swap before \"324.0 -> 435.0\"\n// so the freshly written 435.0 values are not caught by the later search.\nawait replaceAll(context, \"435.0\", \"546.0\");\nawait replaceAll(context, \"324.0\", \"435.0\");\nawait replaceAll(context, \"Bill No. 00017\", \"Bill No. 00018\");\nawait replaceAll(context, \"2025-06-19\", \"2025-06-18\");\n", "ps1": "# Bill revision: 00017 -> 00018, due date 2025-06-19 -> 2025-06-18,\n# meter reading Previous 324.0 -> 435.0, Present 435.0 -> 546.0.\n# The document contains two identical copies of the bill, so Find/Replace\n# (ReplaceAll) catches both occurrences of each target string.\n#\n# Order matters: do the \"435.0 -> 546.0\" swap before \"324.0 -> 435.0\" so the\n# freshly written 435.0 values are not re-matched by the later find.\n\n$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\"435.0\", $false, $false, $false, $false, $false, $true, 1, $false, \"546.0\", 2)\n$d.Content.Find.Execute(\"324.0\", $false, $false, $false, $false, $false, $true, 1, $false, \"435.0\", 2)\n$d.Content.Find.Execute(\"Bill No. 00017\", $false, $false, $false, $false, $false, $true, 1, $false, \"Bill No. 00018\", 2)\n$d.Content.Find.Execute(\"2025-06-19\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-06-18\", 2)\n"}
